$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell without letting Excel's
# "looks like a date/number" auto-detection coerce it (which would also
# stamp the cell with an extra NumberFormat style). We do this by writing
# the string as a formula result (a quoted text formula is always a
# string) into a scratch cell, copying it, and pasting *values only* into
# the destination - this bakes in a plain shared-string cell, same as the
# other plain-text columns.
function Set-TextValue($range, [string]$text) {
    $helper = $ws.Range("Z100")
    $escaped = $text.Replace('"', '""')
    $helper.Formula = "=""" + $escaped + """"
    $helper.Copy()
    $range.PasteSpecial(-4163)
    $helper.ClearContents()
}

# ---- Row 2: update existing transaction ----
$ws.Range("A2").Value = "Food"
$ws.Range("B2").Value = "asdwa"
Set-TextValue $ws.Range("C2") "2023-03-20"
Set-TextValue $ws.Range("D2") "432.0"
$ws.Range("E2").Value = "Savings"

# ---- Row 3: new transaction ----
$ws.Range("A3").Value = "Food"
$ws.Range("B3").Value = "sdaw"
Set-TextValue $ws.Range("C3") "2023-03-23"
Set-TextValue $ws.Range("D3") "432.0"
$ws.Range("E3").Value = "Savings"

# ---- Row 4: new transaction ----
$ws.Range("A4").Value = "Food"
$ws.Range("B4").Value = "dwads"
Set-TextValue $ws.Range("C4") "2023-03-01"
Set-TextValue $ws.Range("D4") "43.0"
$ws.Range("E4").Value = "Savings"

# ---- Row 5: new transaction ----
$ws.Range("A5").Value = "Food"
$ws.Range("B5").Value = "dwas"
Set-TextValue $ws.Range("C5") "2023-03-02"
Set-TextValue $ws.Range("D5") "432.0"
$ws.Range("E5").Value = "Savings"
